$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8-9 get relabelled as new "line7"/"line8" entries, with new C/D values ---
$ws.Cells.Item(8,2).Value = "line7"
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = 11

$ws.Cells.Item(9,2).Value = "line8"
$ws.Cells.Item(9,3).Value = 16

# --- Rows 10-15 (formerly extr3..extr8) shift labels up to extr1..extr6, with new C/D/E values ---
$ws.Cells.Item(10,2).Value = "extr1"
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 12
$ws.Cells.Item(10,5).Value = $true

$ws.Cells.Item(11,2).Value = "extr2"
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 9

$ws.Cells.Item(12,2).Value = "extr3"
$ws.Cells.Item(12,3).Value = 10

$ws.Cells.Item(13,2).Value = "extr4"
$ws.Cells.Item(13,4).Value = 8
$ws.Cells.Item(13,5).Value = $false

$ws.Cells.Item(14,2).Value = "extr5"
$ws.Cells.Item(14,3).Value = 9
$ws.Cells.Item(14,4).Value = 11

$ws.Cells.Item(15,2).Value = "extr6"
$ws.Cells.Item(15,3).Value = 7
$ws.Cells.Item(15,4).Value = 11

# --- New rows 16-17 (extr7, extr8), copying column-A number formatting/style from row 15 ---
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "extr7"
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 7
$ws.Cells.Item(16,5).Value = $true

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "extr8"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = $false
